$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Search")

$ws.Range("C2").Value = "22-09-2020 01:52:22"
$ws.Range("C3").Value = "22-09-2020 01:52:28"
$ws.Range("C4").Value = "22-09-2020 01:52:32"
$ws.Range("C5").Value = "22-09-2020 01:52:38"
